$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Start row 45 as a copy of row 44 (same style/format pattern: date style,
# hyperlink style on G, etc.), then overwrite the values that differ.
$ws.Range("A44:K44").Copy($ws.Range("A45:K45"))

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "fhs"
$ws.Range("C45").Value = 35246
$ws.Range("D45").Value = "Male"
$ws.Range("E45").Value = "Viet Nam"
$ws.Range("F45").Value = 9
$ws.Range("G45").Value = "abcdefs@gmail.com"
$ws.Range("H45").Value = 4
$ws.Range("I45").Value = 45527
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = "Probation"

# Turn G45 into a mailto hyperlink, matching the style used by G42:G44
$ws.Hyperlinks.Add($ws.Range("G45"), "mailto:abcdefs@gmail.com", "", "", "abcdefs@gmail.com")

# Hyperlinks.Add() re-styles the cell with a fresh (duplicate) Hyperlink-like
# style; re-paste the formatting from the row above so G45 reuses the same
# style index as G42:G44.
$ws.Range("G44").Copy()
$ws.Range("G45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep the view roughly where the author left it: scrolled down to row 28,
# with G56 as the active/selected cell.
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G56").Select()
